# "added colors to rows"
# Applies alternating/status background colors to several DTR rows, marks
# row 14 (a holiday/no-work row) with its "late" hour flag, flips B19 from a
# blank-string placeholder to an explicit boolean FALSE, and cleans up the
# extra redundant 3rd argument that had crept into several FLOOR() calls.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Color the "half day" rows (5, 8, 10, 12) with an orange fill, and the
#    "holiday" row (14) with a red fill - while keeping the same font,
#    thin black border and centered/wrapped alignment the rows already had.
# ---------------------------------------------------------------------
$orangeRows = @(5, 8, 10, 12)
foreach ($r in $orangeRows) {
    $rng = $ws.Range("A" + $r + ":J" + $r)
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 11
    $rng.Borders.LineStyle = 1
    $rng.Borders.Color = 0
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
    $rng.WrapText = $true
    $rng.Interior.Color = 6737151
}

$redRange = $ws.Range("A14:J14")
$redRange.Font.Name = "Arial"
$redRange.Font.Size = 11
$redRange.Borders.LineStyle = 1
$redRange.Borders.Color = 0
$redRange.HorizontalAlignment = -4108
$redRange.VerticalAlignment = -4108
$redRange.WrapText = $true
$redRange.Interior.Color = 6184671

# Row 14 is a holiday, so the "no of hours late" column now reflects a
# full day (1) instead of 0.
$ws.Range("I14").Value = 1

# ---------------------------------------------------------------------
# 2. B19 is part of the merged range A19:G19, so a direct .Value assignment
#    is ignored by the merge. Stage the boolean in a scratch cell, copy it
#    across with PasteSpecial (which is allowed to write into merged
#    cells), then clean the scratch cell back up.
# ---------------------------------------------------------------------
$scratch = $ws.Range("Z1")
$scratch.Value = $false
$scratch.Copy()
$ws.Range("B19").PasteSpecial(-4163)
$scratch.ClearContents()

# ---------------------------------------------------------------------
# 3. Drop the redundant third argument Excel was carrying on several
#    FLOOR() calls: FLOOR(x,1,1) -> FLOOR(x,1)
# ---------------------------------------------------------------------
$ws.Range("B22").Formula = '=FLOOR(F17/8,1)&"."&FLOOR(MOD(F17,8),1)&"."&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60'
$ws.Range("B23").Formula = '=FLOOR(H19,1)&"."&(H19-FLOOR(H19,1))*8&".0"'
$ws.Range("B24").Formula = '=FLOOR(I19,1)&"."&(I19-FLOOR(I19,1))*8&".0"'
$ws.Range("B27").Formula = '=FLOOR(K27/8,1)&"."&FLOOR(MOD(K27,8),1)&"."&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60'

Write-Host "done"
